# Apply crypto price/volume updates for Wed Jun 26 14:49:58 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.682.13"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.376.40"
$ws.Range("E3").Value = "  -0.78%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.69"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.89"
$ws.Range("E6").Value = "  -0.38%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.374.87"
$ws.Range("E8").Value = "  -0.74%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.472"
$ws.Range("E9").Value = "  -1.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.46"
$ws.Range("E10").Value = "  -3.07%  "

$ws.Range("E11").Value = "  -0.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.392"
$ws.Range("E12").Value = "  -0.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.949.80"
$ws.Range("E13").Value = "  -0.83%  "

$ws.Range("E14").Value = "  +1.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000175"
$ws.Range("E15").Value = "  -2.41%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.94"
$ws.Range("E16").Value = "  +1.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.376.73"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.702.62"
$ws.Range("E18").Value = "  -0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.91"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.96"
$ws.Range("E20").Value = "  -0.41%  "

$ws.Range("E21").Value = "  -1.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.03"
$ws.Range("E22").Value = "  -3.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.555"
$ws.Range("E23").Value = "  -2.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.522.29"
$ws.Range("E24").Value = "  -0.37%  "

$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.22"
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.79"
$ws.Range("E28").Value = "  +8.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.63"
$ws.Range("E29").Value = "  -2.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.56%  "

$ws.Range("E31").Value = "  +3.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.21"
$ws.Range("E32").Value = "  -0.86%  "

$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("E34").Value = "  -0.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.62"
$ws.Range("E35").Value = "  -0.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.21"
$ws.Range("E36").Value = "  -6.03%  "

$ws.Range("E37").Value = "  -0.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.82"
$ws.Range("E38").Value = "  -3.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.89"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0769"
$ws.Range("E40").Value = "  -3.24%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.74"
$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("E43").Value = "  +0.44%  "

$ws.Range("E44").Value = "  -1.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.53"
$ws.Range("E45").Value = "  +0.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.38"
$ws.Range("E46").Value = "  -1.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.26"
$ws.Range("E47").Value = "  +3.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.85"
$ws.Range("E48").Value = "  -2.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.04"
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.377.26"
$ws.Range("E50").Value = "  +1.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0262"
$ws.Range("E51").Value = "  -1.40%  "
